$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had two new weekly price records inserted before the existing
# row that used to be row 85 ("Fruta, Macroferia Regional de Talca - Ciruela").
# Insert two blank rows at 85, pushing the old rows 85-87 down to 87-89.
$ws.Rows("85:86").Insert()

# New row 85: Black Amber / Especial, Región de O'Higgins, $/caja
$ws.Cells.Item(85, 1).Value = 5
$ws.Cells.Item(85, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(85, 3).Value = "Maule"
$ws.Cells.Item(85, 4).Value = 44610
$ws.Cells.Item(85, 5).Value = 7
$ws.Cells.Item(85, 6).Value = "Fruta"
$ws.Cells.Item(85, 7).Value = 100103
$ws.Cells.Item(85, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(85, 9).Value = 100103002
$ws.Cells.Item(85, 10).Value = "Ciruela"
$ws.Cells.Item(85, 11).Value = "Black Amber"
$ws.Cells.Item(85, 12).Value = "Especial"
$ws.Cells.Item(85, 13).Value = 200
$ws.Cells.Item(85, 14).Value = 13000
$ws.Cells.Item(85, 15).Value = 13000
$ws.Cells.Item(85, 16).Value = 13000
$ws.Cells.Item(85, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(85, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(85, 19).Value = 722
$ws.Cells.Item(85, 20).Value = 18

# New row 86: Black Amber / Primera, Región de O'Higgins, $/caja
$ws.Cells.Item(86, 1).Value = 5
$ws.Cells.Item(86, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(86, 3).Value = "Maule"
$ws.Cells.Item(86, 4).Value = 44610
$ws.Cells.Item(86, 5).Value = 7
$ws.Cells.Item(86, 6).Value = "Fruta"
$ws.Cells.Item(86, 7).Value = 100103
$ws.Cells.Item(86, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(86, 9).Value = 100103002
$ws.Cells.Item(86, 10).Value = "Ciruela"
$ws.Cells.Item(86, 11).Value = "Black Amber"
$ws.Cells.Item(86, 12).Value = "Primera"
$ws.Cells.Item(86, 13).Value = 150
$ws.Cells.Item(86, 14).Value = 11000
$ws.Cells.Item(86, 15).Value = 11000
$ws.Cells.Item(86, 16).Value = 11000
$ws.Cells.Item(86, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(86, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(86, 19).Value = 611
$ws.Cells.Item(86, 20).Value = 18
